# Fixed naive component forecaster bug - Presentation state 11.02.
#
# For each data row (2..16) a new first error value was inserted into
# column B; all the previously existing values in that row shift one
# column to the right, and since the sheet is capped at column K the
# value that used to sit in column K (if the row already reached it)
# falls off the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values to be placed into column B for each row, pushing the
# rest of that row's contents one column to the right.
$newValues = @{
    2  = -1.097142175261494
    3  = 0.5056247995153902
    4  = -1.303839698193279
    5  = 1.784934712404416
    6  = -1.285852883620195
    7  = 0.1655615342000891
    8  = -0.1020898895371165
    9  = 0.2502934172212692
    10 = -0.5750606441290271
    11 = 0.3545997876350467
    12 = 0.1319134556777877
    13 = 0.4278546843610848
    14 = -0.7714259786200386
    15 = 0.6110347010110101
    16 = -0.343237405067616
}

$lastCol = 11   # column K is the fixed right-hand edge of the data block

for ($row = 2; $row -le 16; $row++) {

    # Collect the existing values currently in columns B.. (stop at the
    # first empty cell) for this row.
    $existing = @()
    for ($col = 2; $col -le $lastCol; $col++) {
        $cellValue = $ws.Cells.Item($row, $col).Value()
        if ($cellValue -eq $null) {
            break
        }
        $existing += $cellValue
    }

    # Shift the collected values one column to the right, starting from
    # the right-most column and working left so nothing is overwritten
    # prematurely. Anything that would land beyond column K is dropped.
    for ($i = $existing.Length - 1; $i -ge 0; $i--) {
        $destCol = 2 + $i + 1
        if ($destCol -le $lastCol) {
            $ws.Cells.Item($row, $destCol).Value = $existing[$i]
        }
    }

    # Place the new value into column B.
    $ws.Cells.Item($row, 2).Value = $newValues[$row]
}
